$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Samples" sheet: add a new Plate-Reader sample in row 2 and push the
# existing Flow sample down to row 3 (now well A1 / sample index 2).
# ---------------------------------------------------------------------------
$samples = $wb.Worksheets.Item("Samples")

# Move the existing Flow row (was row 2) down to row 3, updating the
# Data Location path along the way (repo renamed eebio-tools -> esm).
$samples.Range("A3").Value = 2
$samples.Range("B3").Value = "A1"
$samples.Range("D3").Value = "Flow"
$samples.Range("E3").Value = "/Users/qr24461/OneDrive - University of Bristol/Code/esm/test/inputs/small.fcs"
$samples.Range("F3").Value = "FSC-H,SSC-H,FL1-H,FL1-H,FL3-H,FL1-A,FL4-H"

# New row 2: Plate Reader sample read from a directory.
$samples.Range("A2").Value = 1
$samples.Range("B2").Value = ""
$samples.Range("D2").Value = "Plate Reader"
$samples.Range("E2").Value = "/Users/qr24461/OneDrive - University of Bristol/Code/esm/test/inputs/pr_folder"
$samples.Range("F2").Value = "700,(558,602)"
$samples.Range("G2").Value = "tecan"

$samples.Range("E5").Select() | Out-Null

# ---------------------------------------------------------------------------
# "ID" sheet: fill in the Current/Target mapping rows.
# ---------------------------------------------------------------------------
$id = $wb.Worksheets.Item("ID")
$id.Range("A2").Value = "558,602"
$id.Range("B2").Value = "flo"
$id.Range("A3").Value = "700"
$id.Range("B3").Value = "OD"

$id.Range("B4").Select() | Out-Null

# ---------------------------------------------------------------------------
# "Transformations" sheet: update the flow_cyt equation to use channel
# names without the "-H" suffix.
# ---------------------------------------------------------------------------
$transformations = $wb.Worksheets.Item("Transformations")
$transformations.Range("B2").Value = 'process_fcs("plate_01",["FSC","SSC"],["FL1"])'

$transformations.Range("B2").Select() | Out-Null

# ---------------------------------------------------------------------------
# Other sheets: only the cursor/selection moved in the source session;
# replicate that for fidelity.
# ---------------------------------------------------------------------------
$groups = $wb.Worksheets.Item("Groups")
$groups.Range("C19").Select() | Out-Null

$views = $wb.Worksheets.Item("Views")
$views.Range("A7").Select() | Out-Null

$samples.Activate() | Out-Null
